$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ACE_landing_page_data")

# Refresh the ACE landing page data table with the latest figures: a new
# 2022 row is added at the top and every other year's row shifts down one
# slot (the old 2016 row falls off the bottom of the A2:O7 range).
$ws.Range("A2").Value = 2022
$ws.Range("B2").Value = 470.15663229742336
$ws.Range("C2").Value = 8921474338.786293
$ws.Range("D2").Value = 18975536.504061323
$ws.Range("E2").Value = 0.8836954720487562
$ws.Range("F2").Value = 133.16679170159205
$ws.Range("G2").Value = 319.46355316139619
$ws.Range("H2").Value = -0.34650568426651729
$ws.Range("I2").Value = 0.034409572293988999
$ws.Range("J2").Value = 0.58288993092918751
$ws.Range("K2").Value = 0.46856790462369879
$ws.Range("L2").Value = -0.023420464129504381
$ws.Range("M2").Value = -0.35179038408891561
$ws.Range("N2").Value = 96.927408574309311
$ws.Range("O2").Value = 93.290786227371314
$ws.Range("A3").Value = 2021
$ws.Range("B3").Value = 719.45022470427932
$ws.Range("C3").Value = 8624702030.7452507
$ws.Range("D3").Value = 11987906.507764762
$ws.Range("E3").Value = 0.60173960582040065
$ws.Range("F3").Value = 136.36041593160243
$ws.Range("G3").Value = 492.83988592544631
$ws.Range("H3").Value = -0.25294214079289401
$ws.Range("I3").Value = -0.049983832135574113
$ws.Range("J3").Value = 0.27167682684274408
$ws.Range("K3").Value = 0.25122106103595776
$ws.Range("L3").Value = -0.082189450562490274
$ws.Range("M3").Value = -0.246553986261564
$ws.Range("N3").Value = 93.703124149707335
$ws.Range("O3").Value = 58.937001496122846
$ws.Range("A4").Value = 2020
$ws.Range("B4").Value = 963.04485099436852
$ws.Range("C4").Value = 9078479211.7096443
$ws.Range("D4").Value = 9426849.8526687324
$ws.Range("E4").Value = 0.48092189666483548
$ws.Range("F4").Value = 148.57141924896419
$ws.Range("G4").Value = 654.11439829654353
$ws.Range("H4").Value = 1.2134441067149506
$ws.Range("I4").Value = -0.04188692214530354
$ws.Range("J4").Value = -0.56713924921435421
$ws.Range("K4").Value = -0.50702940127839669
$ws.Range("L4").Value = 0.0847184440380222
$ws.Range("M4").Value = 1.2196723111014927
$ws.Range("N4").Value = 98.633188907032832
$ws.Range("O4").Value = 46.345895633286567
$ws.Range("A5").Value = 2019
$ws.Range("B5").Value = 435.08885002913263
$ws.Range("C5").Value = 9475373441.3449364
$ws.Range("D5").Value = 21778019.456739664
$ws.Range("E5").Value = 0.97555898447490952
$ws.Range("F5").Value = 136.96772657048723
$ws.Range("G5").Value = 294.68962379043467
$ws.Range("H5").Value = -0.0024139029120445743
$ws.Range("I5").Value = 0.01411651830421401
$ws.Range("J5").Value = 0.016570420602805447
$ws.Range("K5").Value = 0.010149067284991542
$ws.Range("L5").Value = 0.0061435452387501588
$ws.Range("M5").Value = -0.001673078925594762
$ws.Range("N5").Value = 102.94524851688867
$ws.Range("O5").Value = 107.06883345086935
$ws.Range("A6").Value = 2018
$ws.Range("B6").Value = 436.14165363691069
$ws.Range("C6").Value = 9343476090.0941372
$ws.Range("D6").Value = 21423030.825376313
$ws.Range("E6").Value = 0.96575744716267398
$ws.Range("F6").Value = 136.13139717353735
$ws.Range("G6").Value = 295.18348906517315
$ws.Range("H6").Value = -0.036186849785866837
$ws.Range("I6").Value = 0.015122489958370178
$ws.Range("J6").Value = 0.053235774727536489
$ws.Range("K6").Value = 0.0491400158657882
$ws.Range("L6").Value = -0.00030022649526117995
$ws.Range("M6").Value = -0.03087473276982311
$ws.Range("N6").Value = 101.51224899583701
$ws.Range("O6").Value = 105.32357747275366
$ws.Range("A7").Value = 2017
$ws.Range("B7").Value = 452.51681152100053
$ws.Range("C7").Value = 9204284391.8050823
$ws.Range("D7").Value = 20340204.291786686
$ws.Range("E7").Value = 0.92052293550703623
$ws.Range("F7").Value = 136.17227969982335
$ws.Range("G7").Value = 304.58754822152844
$ws.Range("H7").Value = -0.035869632351568481
$ws.Range("I7").Value = 0.0089638823037501147
$ws.Range("J7").Value = 0.046501506600886477
$ws.Range("K7").Value = 0.042139905519349652
$ws.Range("L7").Value = 0.011321105465380787
$ws.Range("M7").Value = -0.038898520130888015
$ws.Range("N7").Value = 100
$ws.Range("O7").Value = 100

# Make the ACE_landing_page_data tab the active/selected sheet and select
# the refreshed data range, matching the saved view state of the workbook.
$ws.Activate()
$ws.Range("A2:O7").Select()
